# Sample Project / Main.xlsx — resave edit.
#
# The only meaningful cell change is B11 on the active sheet: it held the
# shared string "R40" and must now hold the (text) string "1".
#
# A plain `$cell.Value = "1"` would get auto-coerced to the number 1 by
# Excel's General-format type inference, which is not what the source
# file shows (B11 keeps its shared-string cell type `t="s"`, just
# pointing at a new shared string "1"). To force Excel to store it as
# text instead of a number, enter it as a formula that evaluates to the
# text "1", then convert that formula to its resulting value in place
# via copy / paste-special-values. This preserves the cell's existing
# style (s="23") and yields a genuine text cell, matching the target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""   # formula whose result is the text "1"
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues: bake the formula down to a text value
$excel.CutCopyMode = $false
